$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new rows 5-13.
# Each entry is an ordered list of (column, value, isText) covering A..S.

$rows = @(
    @{ rowNum=5;  A="resnet";  B="2025-09-26 16:27:11"; C="uliege"; D="/home/labsig/Documents/Axelle/Main research/Data/our/validation"; E="all";
       F=0.8184; G=0.945; H=0.9458; I=0.9293; J=0.9664; K=0.9671999999999999; L=0.843; M=0.9452; N=0.9461000000000001;
       O=399.3357; P=313.5957; Q=76.73990000000001; R=2.4342 },
    @{ rowNum=6;  A="hoptim";  B="2025-09-26 17:15:11"; C="uliege"; D="/home/labsig/Documents/Axelle/Main research/Data/our/validation"; E="all";
       F=0.7585; G=0.906; H=0.9173; I=0.8879; J=0.9405; K=0.9506; L=0.7842; M=0.9192; N=0.9314;
       O=1772.5035; P=1624.6798; Q=139.8535; R=2.4627 },
    @{ rowNum=7;  A="uni2";    B="2025-09-26 18:08:42"; C="uliege"; D="/home/labsig/Documents/Axelle/Main research/Data/our/validation"; E="all";
       F=0.7889; G=0.9185; H=0.9304; I=0.9016999999999999; J=0.9415; K=0.9517; L=0.7988; M=0.9177; N=0.9308999999999999;
       O=1629.0576; P=1483.1252; Q=138.4627; R=2.3626 },
    @{ rowNum=8;  A="hoptim";  B="2025-10-16 13:29:44"; C="uliege"; D="/home/labsig/Documents/Axelle/Main research/Data/uliege/sub_fold_test/val"; E="all";
       F=0.7399; G=1; H=1; I=0.8786; J=1; K=1; L=0.7457; M=1; N=1;
       O=2.9859; P=2.8889; Q=0.0835; R=0.0041 },
    @{ rowNum=9;  A="hoptim";  B="2025-10-16 13:47:10"; C="uliege"; D="/home/labsig/Documents/Axelle/Main research/Data/uliege/sub_fold_test/val"; E="all";
       F="cells_no_aug_1"; G=0.3654; H=1; I=1; J=0.5962; K=1; L=1; M=0.2692; N=1;
       O=1; P=0.8982; Q=0.8682; R=0.0258; S=0.0012 },
    @{ rowNum=10; A="hoptim";  B="2025-10-16 13:48:04"; C="uliege"; D="/home/labsig/Documents/Axelle/Main research/Data/uliege/sub_fold_test/val"; E="all";
       F="cells_no_aug_0"; G=0.9008; H=1; I=1; J=1; K=1; L=1; M=0.9504; N=1;
       O=1; P=2.1195; Q=2.054; R=0.0553; S=0.0031 },
    @{ rowNum=11; A="hoptim1"; B="2025-10-16 14:22:33"; C="uliege"; D="/home/labsig/Documents/Axelle/Main research/Data/uliege/sub_fold_test/val"; E="all";
       F="cells_no_aug_0"; G=0.8843; H=1; I=1; J=0.9917; K=1; L=1; M=0.9174; N=1;
       O=1; P=2.1013; Q=2.0305; R=0.0604; S=0.0031 },
    @{ rowNum=12; A="hoptim1"; B="2025-10-16 14:23:14"; C="uliege"; D="/home/labsig/Documents/Axelle/Main research/Data/uliege/sub_fold_test/val"; E="all";
       F="cells_no_aug_1"; G=0.3846; H=1; I=1; J=0.5385; K=1; L=1; M=0.3654; N=1;
       O=1; P=0.9251; Q=0.8939; R=0.0266; S=0.0013 },
    @{ rowNum=13; A="uni2";    B="2025-10-16 14:25:33"; C="uliege"; D="/home/labsig/Documents/Axelle/Main research/Data/uliege/sub_fold_test/val"; E="all";
       F="cells_no_aug_1"; G=0.2308; H=1; I=1; J=0.6731; K=1; L=1; M=0.2308; N=1;
       O=1; P=0.773; Q=0.7409; R=0.0272; S=0.0014 }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S")

foreach ($row in $rows) {
    $r = $row["rowNum"]
    foreach ($col in $cols) {
        if ($row.ContainsKey($col)) {
            $val = $row[$col]
            if ($null -ne $val) {
                $ws.Range("$col$r").Value = $val
            }
        }
    }
}
